$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 75; existing rows 75:86 shift down to 76:87
$ws.Rows.Item(75).Insert()

# Populate the newly inserted row 75 with the new record
$ws.Cells.Item(75, 1).Value = 8
$ws.Cells.Item(75, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(75, 3).Value = "Coquimbo"
$ws.Cells.Item(75, 4).Value = 44694
$ws.Cells.Item(75, 5).Value = 4
$ws.Cells.Item(75, 6).Value = 100112052
$ws.Cells.Item(75, 7).Value = "Albahaca"
$ws.Cells.Item(75, 8).Value = "Sin especificar"
$ws.Cells.Item(75, 9).Value = "Primera"
$ws.Cells.Item(75, 10).Value = 1000
$ws.Cells.Item(75, 11).Value = 5000
$ws.Cells.Item(75, 12).Value = 5500
$ws.Cells.Item(75, 13).Value = 5250
$ws.Cells.Item(75, 14).Value = '$/docena de matas'
$ws.Cells.Item(75, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(75, 16).Value = 875
$ws.Cells.Item(75, 17).Value = 6
$ws.Cells.Item(75, 18).Value = "Hortaliza"

Write-Output "Row inserted and populated"
